# Commit by Farheen: Resolve the bug InVan-35, and did validation in Item master module.
# Update the sample/demo data row in the "Download-Format-Company" template
# so it no longer collides with values a real import would treat as
# duplicates/placeholders, and flip the IsBlackListed sample flag to False.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: CompanyName sample value "ABC company name" -> "DEF company"
$ws.Range("A2").Value = "DEF company"

# D2: ContactPersonName sample value "TM_0001" -> "Sunny"
$ws.Range("D2").Value = "Sunny"

# J2: IsBlackListed sample flag TRUE -> FALSE
$ws.Range("J2").Value = $false

# Leave the cursor where the author left it when they saved the file
$ws.Range("J6").Select()
